$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp in A1 (14:52 -> 15:22)
$ws.Range("A1").Value = "Datos actualizados a 24 de Abril de 2020 a las 15:22"

# Row 4 - Estados Unidos: update case counts (D4/F4 unchanged)
$ws.Range("B4").Value = 887622
$ws.Range("C4").Value = 1180
$ws.Range("E4").Value = 751417
$ws.Range("G4").Value = 47
$ws.Range("H4").Value = 50283

# Row 20: update case counts (all columns)
$ws.Range("B20").Value = 22797
$ws.Range("C20").Value = 444
$ws.Range("D20").Value = 1228
$ws.Range("E20").Value = 20715
$ws.Range("F20").Value = 188
$ws.Range("G20").Value = 34
$ws.Range("H20").Value = 854

# Rows 24-26: Arabia Saudita overtakes Austria and Israel in the ranking,
# so the three countries shift down one row each, with Arabia Saudita
# getting freshly updated figures and Austria/Israel keeping their totals.
$ws.Range("A24").Value = "Arabia Saudita"
$ws.Range("B24").Value = 15102
$ws.Range("C24").Value = 1172
$ws.Range("D24").Value = 2049
$ws.Range("E24").Value = 12926
$ws.Range("F24").Value = 93
$ws.Range("G24").Value = 6
$ws.Range("H24").Value = 127

$ws.Range("A25").Value = "Austria"
$ws.Range("B25").Value = 15071
$ws.Range("C25").Value = 69
$ws.Range("D25").Value = 11872
$ws.Range("E25").Value = 2677
$ws.Range("F25").Value = 169
$ws.Range("G25").Value = 0
$ws.Range("H25").Value = 522

$ws.Range("A26").Value = "Israel"
$ws.Range("B26").Value = 14882
$ws.Range("C26").Value = 79
$ws.Range("D26").Value = 5685
$ws.Range("E26").Value = 9004
$ws.Range("F26").Value = 139
$ws.Range("G26").Value = 1
$ws.Range("H26").Value = 193

# Row 64: update case counts (D64/F64 unchanged)
$ws.Range("B64").Value = 2376
$ws.Range("C64").Value = 87
$ws.Range("E64").Value = 1749
$ws.Range("G64").Value = 5
$ws.Range("H64").Value = 25

# Row 80: update case counts (B80/C80/F80 unchanged)
$ws.Range("D80").Value = 188
$ws.Range("E80").Value = 1099
$ws.Range("G80").Value = 1
$ws.Range("H80").Value = 43

# Row 111: update case counts (B111/C111/F111 unchanged)
$ws.Range("D111").Value = 109
$ws.Range("E111").Value = 298
